$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"21.194154"
$ws.Range("H2").Value = [double]"42.38830799999999"
$ws.Range("I2").Value = [double]"0.02540793187283288"
$ws.Range("J2").Value = [double]"0.01757716018274786"
$ws.Range("O2").Value = [double]"0.9919525181111984"
$ws.Range("P2").Value = [double]"0.9919525181111983"
$ws.Range("Q2").Value = [double]"0.8028910712639999"
$ws.Range("R2").Value = [double]"4.817346427583999"
$ws.Range("S2").Value = [double]"0.02520346200125435"
$ws.Range("T2").Value = [double]"0.01743570830452063"

# Row 3
$ws.Range("D3").Value = "Neutrophils"
$ws.Range("G3").Value = [double]"21.194154"
$ws.Range("H3").Value = [double]"42.38830799999999"
$ws.Range("I3").Value = [double]"0.02540793187283288"
$ws.Range("J3").Value = [double]"0.01757716018274786"
$ws.Range("M3").Value = [double]"0.0003073333333333333"
$ws.Range("N3").Value = [double]"0.000922"
$ws.Range("O3").Value = [double]"0.008047481888801606"
$ws.Range("P3").Value = [double]"0.008047481888801606"
$ws.Range("Q3").Value = [double]"0.006513669995999999"
$ws.Range("R3").Value = [double]"0.03908201997599999"
$ws.Range("S3").Value = [double]"0.0002044698715785276"
$ws.Range("T3").Value = [double]"0.0001414518782272281"

# Row 4
$ws.Range("I4").Value = [double]"0.890627797432279"
$ws.Range("J4").Value = [double]"0.9242019895810776"
$ws.Range("O4").Value = [double]"0.9919525181111984"
$ws.Range("P4").Value = [double]"0.9919525181111983"
$ws.Range("S4").Value = [double]"0.8834604863627795"
$ws.Range("T4").Value = [double]"0.9167644908083293"

# Row 5
$ws.Range("D5").Value = "Neutrophils"
$ws.Range("I5").Value = [double]"0.890627797432279"
$ws.Range("J5").Value = [double]"0.9242019895810776"
$ws.Range("M5").Value = [double]"0.0003073333333333333"
$ws.Range("N5").Value = [double]"0.000922"
$ws.Range("O5").Value = [double]"0.008047481888801606"
$ws.Range("P5").Value = [double]"0.008047481888801606"
$ws.Range("Q5").Value = [double]"0.2283245874073333"
$ws.Range("R5").Value = [double]"2.054921286666"
$ws.Range("S5").Value = [double]"0.007167311069499531"
$ws.Range("T5").Value = [double]"0.007437498772748132"

# Row 6
$ws.Range("E6").Value = [double]"2"
$ws.Range("F6").Value = [double]"0.6666666666666666"
$ws.Range("G6").Value = [double]"0.05889633333333333"
$ws.Range("H6").Value = [double]"0.176689"
$ws.Range("I6").Value = [double]"7.060598054034096E-05"
$ws.Range("J6").Value = [double]"7.326762973246153E-05"
$ws.Range("O6").Value = [double]"0.9919525181111984"
$ws.Range("P6").Value = [double]"0.9919525181111983"
$ws.Range("Q6").Value = [double]"0.002231150163555555"
$ws.Range("R6").Value = [double]"0.020080351472"
$ws.Range("S6").Value = [double]"7.003778019070149E-05"
$ws.Range("T6").Value = [double]"7.267800980915412E-05"

# Row 7
$ws.Range("D7").Value = "Neutrophils"
$ws.Range("E7").Value = [double]"2"
$ws.Range("F7").Value = [double]"0.6666666666666666"
$ws.Range("G7").Value = [double]"0.05889633333333333"
$ws.Range("H7").Value = [double]"0.176689"
$ws.Range("I7").Value = [double]"7.060598054034096E-05"
$ws.Range("J7").Value = [double]"7.326762973246153E-05"
$ws.Range("M7").Value = [double]"0.0003073333333333333"
$ws.Range("N7").Value = [double]"0.000922"
$ws.Range("O7").Value = [double]"0.008047481888801606"
$ws.Range("P7").Value = [double]"0.008047481888801606"
$ws.Range("Q7").Value = [double]"1.810080644444444E-05"
$ws.Range("R7").Value = [double]"0.000162907258"
$ws.Range("S7").Value = [double]"5.682003496394725E-07"
$ws.Range("T7").Value = [double]"5.896199233074062E-07"

# Row 8
$ws.Range("G8").Value = [double]"69.71480750000001"
$ws.Range("H8").Value = [double]"139.429615"
$ws.Range("I8").Value = [double]"0.08357536137029385"
$ws.Range("J8").Value = [double]"0.05781728011115387"
$ws.Range("O8").Value = [double]"0.9919525181111984"
$ws.Range("P8").Value = [double]"0.9919525181111983"
$ws.Range("Q8").Value = [double]"2.640982814253334"
$ws.Range("R8").Value = [double]"15.84589688552"
$ws.Range("S8").Value = [double]"0.08290279016331636"
$ws.Range("T8").Value = [double]"0.05735199659659958"

# Row 9
$ws.Range("D9").Value = "Neutrophils"
$ws.Range("G9").Value = [double]"69.71480750000001"
$ws.Range("H9").Value = [double]"139.429615"
$ws.Range("I9").Value = [double]"0.08357536137029385"
$ws.Range("J9").Value = [double]"0.05781728011115387"
$ws.Range("M9").Value = [double]"0.0003073333333333333"
$ws.Range("N9").Value = [double]"0.000922"
$ws.Range("O9").Value = [double]"0.008047481888801606"
$ws.Range("P9").Value = [double]"0.008047481888801606"
$ws.Range("Q9").Value = [double]"0.02142568417166667"
$ws.Range("R9").Value = [double]"0.12855410503"
$ws.Range("S9").Value = [double]"0.0006725712069774891"
$ws.Range("T9").Value = [double]"0.00046528351455428"

# Row 10
$ws.Range("G10").Value = [double]"0.1138996666666667"
$ws.Range("H10").Value = [double]"0.341699"
$ws.Range("I10").Value = [double]"0.0001365449628706596"
$ws.Range("J10").Value = [double]"0.000141692328395952"
$ws.Range("O10").Value = [double]"0.9919525181111984"
$ws.Range("P10").Value = [double]"0.9919525181111983"
$ws.Range("Q10").Value = [double]"0.004314823105777778"
$ws.Range("R10").Value = [double]"0.038833407952"
$ws.Range("S10").Value = [double]"0.0001354461197549508"
$ws.Range("T10").Value = [double]"0.0001405520619494035"

# Row 11
$ws.Range("D11").Value = "Neutrophils"
$ws.Range("G11").Value = [double]"0.1138996666666667"
$ws.Range("H11").Value = [double]"0.341699"
$ws.Range("I11").Value = [double]"0.0001365449628706596"
$ws.Range("J11").Value = [double]"0.000141692328395952"
$ws.Range("M11").Value = [double]"0.0003073333333333333"
$ws.Range("N11").Value = [double]"0.000922"
$ws.Range("O11").Value = [double]"0.008047481888801606"
$ws.Range("P11").Value = [double]"0.008047481888801606"
$ws.Range("Q11").Value = [double]"3.500516422222222E-05"
$ws.Range("R11").Value = [double]"0.000315046478"
$ws.Range("S11").Value = [double]"1.09884311570872E-06"
$ws.Range("T11").Value = [double]"1.140266446548554E-06"

# Row 12
$ws.Range("G12").Value = [double]"0.1516146666666666"
$ws.Range("H12").Value = [double]"0.454844"
$ws.Range("I12").Value = [double]"0.0001817583811832703"
$ws.Range("J12").Value = [double]"0.0001886101668922895"
$ws.Range("O12").Value = [double]"0.9919525181111984"
$ws.Range("P12").Value = [double]"0.9919525181111983"
$ws.Range("Q12").Value = [double]"0.005743567879111111"
$ws.Range("R12").Value = [double]"0.051692110912"
$ws.Range("S12").Value = [double]"0.00018029568390256"
$ws.Range("T12").Value = [double]"0.0001870923299901799"

# Row 13
$ws.Range("D13").Value = "Neutrophils"
$ws.Range("G13").Value = [double]"0.1516146666666666"
$ws.Range("H13").Value = [double]"0.454844"
$ws.Range("I13").Value = [double]"0.0001817583811832703"
$ws.Range("J13").Value = [double]"0.0001886101668922895"
$ws.Range("M13").Value = [double]"0.0003073333333333333"
$ws.Range("N13").Value = [double]"0.000922"
$ws.Range("O13").Value = [double]"0.008047481888801606"
$ws.Range("P13").Value = [double]"0.008047481888801606"
$ws.Range("Q13").Value = [double]"4.659624088888889E-05"
$ws.Range("R13").Value = [double]"0.0004193661679999999"
$ws.Range("S13").Value = [double]"1.462697280710266E-06"
$ws.Range("T13").Value = [double]"1.517836902109548E-06"
